# Excel COM-interop script applying the commit's changes:
# 1. Rename header "Requested quantity" -> "Weekly_PO_Qty" on "Weekly Quantity" sheet (B1)
# 2. Rename header "Requested quantity" -> "Monthly_PO_Qty" on "Monthly Trend" sheet (B1)
# 3. Add a new "PO Forecast" worksheet (sheetId 3) with forecast data

$wb = $excel.ActiveWorkbook

# --- 1. Update "Weekly Quantity" header ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Update "Monthly Trend" header ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the new "PO Forecast" worksheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the sheetPr / pageMargins metadata used by the other sheets in the workbook
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Copy the header formatting (bold, centered, bordered) from the Weekly Quantity sheet
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows: row, ds(date-serial), PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(2, 45067.99999999999, 0, -515.7270505407473, 338.1890730377376),
    @(3, 45074.99999999999, 0, -458.9684968094855, 340.1596254784321),
    @(4, 45221.99999999999, 82, -341.9266694413762, 486.380750045059),
    @(5, 45228.99999999999, 91, -325.3594239348914, 523.3937030244689),
    @(6, 45235.99999999999, 99, -319.0085282766021, 519.9147780455847),
    @(7, 45242.99999999999, 107, -308.6756899475969, 524.0234205966356),
    @(8, 45249.99999999999, 116, -307.328046497262, 512.9433483481139),
    @(9, 45256.99999999999, 124, -320.1939176717243, 534.7431169356089),
    @(10, 45263.99999999999, 132, -263.688708895115, 563.8342197997064),
    @(11, 45270.99999999999, 140, -290.453891899103, 580.1838671960778),
    @(12, 45277.99999999999, 149, -253.6020124581906, 592.7459756101374),
    @(13, 45298.99999999999, 173, -220.0611728808311, 568.6275149486105),
    @(14, 45305.99999999999, 182, -247.5776838498852, 583.98373820581),
    @(15, 45312.99999999999, 190, -225.2452090804507, 632.1997047241188),
    @(16, 45326.99999999999, 206, -216.6634407210675, 624.2845709900257),
    @(17, 45354.99999999999, 239, -187.3619534810607, 668.7637738540077),
    @(18, 45361.99999999999, 248, -170.9412398022154, 678.5083634003256),
    @(19, 45368.99999999999, 256, -171.5541636646221, 672.8928996725297),
    @(20, 45375.99999999999, 264, -150.5383803138391, 687.523539696974),
    @(21, 45382.99999999999, 273, -162.0378978646252, 670.4812346227976),
    @(22, 45389.99999999999, 281, -152.2283526417381, 712.0738902954122),
    @(23, 45396.99999999999, 289, -102.9817118848478, 706.6359561641843),
    @(24, 45417.99999999999, 314, -121.8780089849716, 729.8613886163072),
    @(25, 45424.99999999999, 322, -106.5760469480262, 714.5477092312924),
    @(26, 45459.99999999999, 363, -91.30016702122958, 783.6282861203081),
    @(27, 45466.99999999999, 372, -61.05240080710337, 780.0277048996452),
    @(28, 45494.99999999999, 405, -25.12621184026153, 818.5765141803724),
    @(29, 45501.99999999999, 413, -5.335039544663375, 850.9859233237447),
    @(30, 45515.99999999999, 430, 18.88172944013869, 861.2384469303162),
    @(31, 45522.99999999999, 438, -7.446109467718404, 865.9935344999835),
    @(32, 45543.99999999999, 463, 66.3921387800163, 864.2318804557813),
    @(33, 45550.99999999999, 471, 73.56044761400678, 896.9451335965547),
    @(34, 45557.99999999999, 479, 34.18815236842609, 872.2822414696055),
    @(35, 45564.99999999999, 487, 54.290270676591, 876.4841749575692),
    @(36, 45571.99999999999, 496, 83.65766662209815, 892.810385971933),
    @(37, 45578.99999999999, 504, 40.50501145187719, 930.2027533174114),
    @(38, 45585.99999999999, 512, 114.6606003074603, 950.0679035605334),
    @(39, 45613.99999999999, 545, 125.4483479592025, 971.9017959553229),
    @(40, 45620.99999999999, 554, 146.7650464605377, 992.9685625939452),
    @(41, 45627.99999999999, 562, 135.6010240304861, 977.5039815703088),
    @(42, 45634.99999999999, 570, 152.8838911517851, 983.3764802455919),
    @(43, 45641.99999999999, 578, 167.8375251778789, 993.4002369506148),
    @(44, 45648.99999999999, 587, 167.0084078536853, 1016.360181268529),
    @(45, 45655.99999999999, 595, 182.8510627772892, 1000.296298611456),
    @(46, 45662.99999999999, 603, 180.7947332774344, 1031.81027634864),
    @(47, 45669.99999999999, 611, 193.4521947089331, 1035.919883254517)
)

foreach ($row in $data) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
    $wsForecast.Cells.Item($r, 3).Value = $row[3]
    $wsForecast.Cells.Item($r, 4).Value = $row[4]
}

# Copy the date-number-format style (used in column A of the other sheets) onto the
# "ds" column of the new sheet
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A47").PasteSpecial(-4122)
